$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update F33 text (do this first so the shared-string table append order
# matches what Excel produced in the real edit)
$ws.Range("F33").Value = "valutazione didattica, revisione progetti, presenta Master DS e possibilità tesi su ANN bayesiane e  gerarchiche"

# Update F30: "riprendere lin reg" -> "lin reg"
$ws.Range("F30").Value = "lin reg"

# Add new cell values F31 and F32: "lin reg; reg gerarchica"
$ws.Range("F31").Value = "lin reg; reg gerarchica"
$ws.Range("F32").Value = "lin reg; reg gerarchica"

# Update sheet view (zoom, top-left cell, selection) to match final state
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 150
$win.ScrollColumn = 4
$win.ScrollRow = 16
$ws.Range("F32").Select()
